$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "89.073.11"
Set-TextValue "E2" "  -1.56%  "
Set-TextValue "D3" "3.099.89"
Set-TextValue "E3" "  -3.09%  "
Set-TextValue "E4" "  -0.20%  "
Set-TextValue "D5" "212.89"
Set-TextValue "E5" "  -3.74%  "
Set-TextValue "D6" "623.30"
Set-TextValue "E6" "  -2.25%  "
Set-TextValue "E7" "  -5.30%  "
Set-TextValue "D8" "0.807"
Set-TextValue "E8" "  +14.62%  "
Set-TextValue "E9" "  -0.09%  "
Set-TextValue "D10" "3.102.76"
Set-TextValue "E10" "  -2.95%  "
Set-TextValue "D11" "0.596"
Set-TextValue "E11" "  +3.86%  "
Set-TextValue "E12" "  +0.03%  "
Set-TextValue "D13" "0.0000242"
Set-TextValue "E13" "  -5.73%  "
Set-TextValue "E14" "  -2.55%  "
Set-TextValue "D15" "88.542.42"
Set-TextValue "E15" "  -1.92%  "
Set-TextValue "D16" "32.36"
Set-TextValue "E16" "  -3.01%  "
Set-TextValue "D17" "3.668.83"
Set-TextValue "E17" "  -3.26%  "
Set-TextValue "D18" "3.080.37"
Set-TextValue "E18" "  -3.47%  "
Set-TextValue "D19" "3.39"
Set-TextValue "E19" "  +1.36%  "
Set-TextValue "E20" "  -5.39%  "
Set-TextValue "D21" "13.47"
Set-TextValue "E21" "  +0.01%  "
Set-TextValue "D22" "424.19"
Set-TextValue "E22" "  -3.07%  "
Set-TextValue "E23" "  -3.81%  "
Set-TextValue "E24" "  -2.33%  "
Set-TextValue "D25" "5.65"
Set-TextValue "E25" "  +6.14%  "
Set-TextValue "D26" "11.96"
Set-TextValue "E26" "  +0.85%  "
Set-TextValue "D27" "82.59"
Set-TextValue "E27" "  +1.70%  "
Set-TextValue "E28" "  -4.30%  "
Set-TextValue "E29" "  +0.06%  "
Set-TextValue "E30" "  +7.66%  "
Set-TextValue "E31" "  +7.15%  "
Set-TextValue "E32" "  -3.54%  "
Set-TextValue "D33" "512.24"
Set-TextValue "E33" "  -5.02%  "
Set-TextValue "E34" "  -11.39%  "
Set-TextValue "D35" "6.81"
Set-TextValue "E35" "  -3.84%  "
Set-TextValue "E36" "  -2.20%  "
Set-TextValue "E37" "  -5.58%  "
Set-TextValue "D38" "22.34"
Set-TextValue "E38" "  -0.70%  "
Set-TextValue "E39" "  -0.41%  "
Set-TextValue "E40" "  +2.47%  "
Set-TextValue "E42" "  +0.04%  "
Set-TextValue "D43" "0.365"
Set-TextValue "E43" "  -2.36%  "
Set-TextValue "D44" "1.84"
Set-TextValue "E44" "  -5.44%  "
Set-TextValue "E45" "  +0.17%  "
Set-TextValue "E46" "  +5.45%  "
Set-TextValue "D47" "0.0695"
Set-TextValue "E47" "  +14.11%  "
Set-TextValue "E48" "  -3.05%  "
Set-TextValue "D49" "163.35"
Set-TextValue "E49" "  -5.66%  "
Set-TextValue "D50" "1.22"
Set-TextValue "E50" "  -0.81%  "
Set-TextValue "E51" "  -5.61%  "
